# Fruta / hortaliza, semanal
# Insert a new weekly record for "Terminal Hortofrutícola Agro Chillán" / Pera
# as row 183, pushing the existing rows 183:208 down to 184:209.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 183 (shifts 183:208 -> 184:209)
$ws.Rows("183:183").Insert()

# Populate the newly inserted row 183 with the new weekly observation
$ws.Range("A183").Value = 7
$ws.Range("B183").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C183").Value = "Ñuble"
$ws.Range("D183").Value = 44769
$ws.Range("E183").Value = 16
$ws.Range("F183").Value = "Fruta"
$ws.Range("G183").Value = 100104
$ws.Range("H183").Value = "Frutos de pepita"
$ws.Range("I183").Value = 100104005
$ws.Range("J183").Value = "Pera"
$ws.Range("K183").Value = "Packham's Triumph"
$ws.Range("L183").Value = "Primera"
$ws.Range("M183").Value = 120
$ws.Range("N183").Value = 8000
$ws.Range("O183").Value = 9000
$ws.Range("P183").Value = 8500
$ws.Range("Q183").Value = "`$/caja 16 kilos empedrada"
$ws.Range("R183").Value = "Provincia de Curicó"
$ws.Range("S183").Value = 531
$ws.Range("T183").Value = 16
